$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "232.45")
# are preserved as text (matching the original inline-string cell type)
# instead of being auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '37.724.67'
$ws.Range("E2").Value = '  +1.17%  '
$ws.Range("D3").Value = '2.092.93'
$ws.Range("E3").Value = '  +1.58%  '
$ws.Range("D5").Value = '232.45'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").Value = '0.624'
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '57.78'
$ws.Range("E8").Value = '  +1.74%  '
$ws.Range("E9").Value = '  +1.63%  '
$ws.Range("D10").Value = '0.0777'
$ws.Range("E10").Value = '  +2.25%  '
$ws.Range("E11").Value = '  +2.96%  '
$ws.Range("D12").Value = '2.389.28'
$ws.Range("E12").Value = '  +1.10%  '
$ws.Range("D13").Value = '14.45'
$ws.Range("E13").Value = '  -1.07%  '
$ws.Range("D14").Value = '21.07'
$ws.Range("E14").Value = '  +1.99%  '
$ws.Range("E15").Value = '  -1.61%  '
$ws.Range("D16").Value = '5.24'
$ws.Range("E16").Value = '  +1.68%  '
$ws.Range("D17").Value = '2.078.55'
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").Value = '37.670.47'
$ws.Range("E18").Value = '  +1.24%  '
$ws.Range("D19").Value = '6.16'
$ws.Range("E19").Value = '  -2.84%  '
$ws.Range("D20").Value = '70.54'
$ws.Range("E20").Value = '  +1.73%  '
$ws.Range("D21").Value = '0.0₃0822'
$ws.Range("E21").Value = '  +1.51%  '
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").Value = '2.40'
$ws.Range("E24").Value = '  -1.56%  '
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("D26").Value = '167.78'
$ws.Range("E26").Value = '  +0.70%  '
$ws.Range("D27").Value = '0.140'
$ws.Range("E27").Value = '  +9.41%  '
$ws.Range("D28").Value = '8.94'
$ws.Range("E28").Value = '  +2.05%  '
$ws.Range("E29").Value = '  -1.67%  '
$ws.Range("D30").Value = '19.47'
$ws.Range("E30").Value = '  +2.12%  '
$ws.Range("E31").Value = '  +1.18%  '
$ws.Range("E32").Value = '  +3.97%  '
$ws.Range("D33").Value = '0.0626'
$ws.Range("E33").Value = '  +1.24%  '
$ws.Range("D34").Value = '4.58'
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("D35").Value = '2.49'
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("D36").Value = '1.82'
$ws.Range("E36").Value = '  +4.22%  '
$ws.Range("E37").Value = '  +5.14%  '
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("D39").Value = '5.40'
$ws.Range("E39").Value = '  -4.92%  '
$ws.Range("D40").Value = '0.0989'
$ws.Range("E40").Value = '  +5.49%  '
$ws.Range("E41").Value = '  -0.27%  '
$ws.Range("D42").Value = '97.39'
$ws.Range("E42").Value = '  +1.56%  '
$ws.Range("D43").Value = '1.459.55'
$ws.Range("E43").Value = '  -1.03%  '
$ws.Range("D44").Value = '0.0214'
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("E45").Value = '  -0.66%  '
$ws.Range("E46").Value = '  +3.62%  '
$ws.Range("D47").Value = '15.60'
$ws.Range("E47").Value = '  +3.46%  '
$ws.Range("D48").Value = '4.04'
$ws.Range("E48").Value = '  -6.20%  '
$ws.Range("D49").Value = '7.36'
$ws.Range("E49").Value = '  +2.92%  '
$ws.Range("E50").Value = '  +2.01%  '
$ws.Range("D51").Value = '2.285.13'
$ws.Range("E51").Value = '  +1.54%  '

# Restore the original (default) cell style for column D so no stray
# number-format styling is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"
